$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source row for "IDYLIC" (a duplicate/typo of "IDYLLIC26") is removed;
# all rows below it shift up by one.
$ws.Rows.Item(20).Delete()

# Keep the hidden AutoFilter database defined name in sync with the new
# last-row of the data range (autoFilter ref itself stays at its old extent).
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Sheet1!_FilterDatabase") {
        $n.RefersTo = "=Sheet1!`$A`$1:`$E`$34"
    }
}

# Reflect where the cursor ended up after the deletion.
$ws.Range("D20").Select() | Out-Null
